# Practice: Leetcode Problem#14 Longest Common Prefix
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 was missing the "Level" (column C) value - fill it in like the other rows.
$ws.Range("C28").Value = "Easy"

# Copy the formatting of row 28 down into the new row 29 so the new row
# matches the existing style pattern (wrap-text problem column, text-format
# date column, etc.) before writing values.
$ws.Range("A28:G28").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(29).RowHeight = $ws.Rows.Item(28).RowHeight

# New data row for LeetCode Problem #14.
$ws.Range("A29").Value = "LeetCode"
$ws.Range("B29").Value = "Stephan"
$ws.Range("C29").Value = "Easy"
$ws.Range("D29").Value = "14. Longest Common Prefix"
$ws.Range("E29").Value = "2020/12/17"
$ws.Range("F29").Value = "Sring"
$ws.Range("G29").Value = "Completed"

[void]$ws.Range("E32").Select()
